# Weekly update: two new rows of "Apio" price data (week of 2023-04-25,
# serial 45041) are inserted at row 415, pushing the existing rows 415-438
# down to 417-440 (dimension grows from A1:R438 to A1:R440).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at row 415 (everything below shifts down by 2).
$ws.Rows.Item(415).Resize(2, 1).EntireRow.Insert()

# --- Row 415: Apio, Americana (o), Primera -----------------------------
$ws.Cells.Item(415, 1).Value = 9
$ws.Cells.Item(415, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(415, 3).Value = "Metropolitana"
$ws.Cells.Item(415, 4).Value = 45041
$ws.Cells.Item(415, 5).Value = 13
$ws.Cells.Item(415, 6).Value = 100112017
$ws.Cells.Item(415, 7).Value = "Apio"
$ws.Cells.Item(415, 8).Value = "Americana (o)"
$ws.Cells.Item(415, 9).Value = "Primera"
$ws.Cells.Item(415, 10).Value = 70
$ws.Cells.Item(415, 11).Value = 8000
$ws.Cells.Item(415, 12).Value = 9000
$ws.Cells.Item(415, 13).Value = 8500
$ws.Cells.Item(415, 14).Value = "$/docena de matas"
$ws.Cells.Item(415, 15).Value = "Región de Coquimbo"
$ws.Cells.Item(415, 16).Value = 1417
$ws.Cells.Item(415, 17).Value = 6
$ws.Cells.Item(415, 18).Value = "Hortaliza"

# --- Row 416: Apio, Americana (o), Segunda ------------------------------
$ws.Cells.Item(416, 1).Value = 9
$ws.Cells.Item(416, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(416, 3).Value = "Metropolitana"
$ws.Cells.Item(416, 4).Value = 45041
$ws.Cells.Item(416, 5).Value = 13
$ws.Cells.Item(416, 6).Value = 100112017
$ws.Cells.Item(416, 7).Value = "Apio"
$ws.Cells.Item(416, 8).Value = "Americana (o)"
$ws.Cells.Item(416, 9).Value = "Segunda"
$ws.Cells.Item(416, 10).Value = 52
$ws.Cells.Item(416, 11).Value = 7000
$ws.Cells.Item(416, 12).Value = 7000
$ws.Cells.Item(416, 13).Value = 7000
$ws.Cells.Item(416, 14).Value = "$/docena de matas"
$ws.Cells.Item(416, 15).Value = "Región de Coquimbo"
$ws.Cells.Item(416, 16).Value = 1167
$ws.Cells.Item(416, 17).Value = 6
$ws.Cells.Item(416, 18).Value = "Hortaliza"
